# Agrego CP03 y CP07
# This script reproduces the commit that:
#  - renames the CP03 test-case row (row 5) to "CP03LoginUserInvalido" and
#    fills in its expected-result text in column D
#  - marks the first block of test cases (rows 2-5) as done/verified with a
#    green highlight
#  - turns the old placeholder rows 8 ("CP06 Eliminar item carrito") and 9
#    ("CP07 SubNewsletter") into two fully specified test cases
#    ("LoginPassInvalido" and "CP07SubNewsletter") complete with
#    hyperlinked e-mail addresses and expected-result text
#  - widens column D so the new, longer, expected-result text fits

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5: CP03 test case gets renamed and documented -------------------
$ws.Range("A5").Value = "CP03LoginUserInvalido"
$ws.Range("D5").Value = "Login was unsuccessful. Please correct the errors and try again."

# --- Row 8: former "CP06 Eliminar item carrito" becomes the
#            "LoginPassInvalido" test case -------------------------------
$ws.Range("A8").Value = "LoginPassInvalido "
$ws.Range("B8").Value = "prueba_1@gmail.com"
$ws.Range("C8").Value = 123
$ws.Hyperlinks.Add($ws.Range("B8"), "mailto:prueba_1@gmail.com")

# --- Row 9: former "CP07 SubNewsletter" becomes the fully specified
#            "CP07SubNewsletter" test case -------------------------------
$ws.Range("A9").Value = "CP07SubNewsletter"
$ws.Range("B9").Value = "prueba_1@gmail.com"
$ws.Range("C9").Value = 123456
$ws.Range("D9").Value = "Thank you for signing up! A verification email has been sent. We appreciate your interest."
$ws.Hyperlinks.Add($ws.Range("B9"), "mailto:prueba_1@gmail.com")

# --- Highlight the first block of (now verified) test cases in green -----
$ws.Range("A2:F5").Interior.Color = 5296274

# D4 and D5 (inside the highlighted block) lose their cell border
$ws.Range("D4:D5").Borders.LineStyle = -4142

# --- Column D needs to be much wider to fit the new expected-result text -
$ws.Columns.Item(4).ColumnWidth = 80.86

# --- Selection, as left by the editor, moved to D15 -----------------------
$ws.Range("D15").Select()
